$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 42-69 (A:H) hold one "war participation" record per player (Name,
# Player Status, Fonte de Dados, and 5 "Guerra" history columns). This
# update rotates that block of records up by one row, with the record
# that used to be on row 42 wrapping around to the bottom (row 69).

$firstRow = 42
$lastRow = 69

# Snapshot the original block before we overwrite anything.
$original = $ws.Range("A" + $firstRow + ":H" + $lastRow).Value()

$rowCount = $lastRow - $firstRow + 1
$colCount = 8

$rotated = New-Object 'object[,]' $rowCount, $colCount

for ($r = 1; $r -le $rowCount; $r++) {
    # source row: the next row, wrapping the very first row back to the end
    $srcRow = $r + 1
    if ($srcRow -gt $rowCount) {
        $srcRow = 1
    }
    for ($c = 1; $c -le $colCount; $c++) {
        $rotated[$r - 1, $c - 1] = $original[$srcRow, $c]
    }
}

$ws.Range("A" + $firstRow + ":H" + $lastRow).Value = $rotated
